$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33: new match data row (24/10/2025 Sporting Cristal 0-1 U. de Deportes)
$ws.Range("A33").Value = "24/10/2025"
$ws.Range("B33").Value = "Sporting Cristal"
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = "U. de Deportes"
$ws.Range("F33").Value = "W"
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 1
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1.32
$ws.Range("L33").Value = 0.8
$ws.Range("M33").Value = 8
$ws.Range("N33").Value = 16
$ws.Range("O33").Value = 3
$ws.Range("P33").Value = 4
